$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.878.71"
$ws.Range("E2").Value = "  -3.42%  "
$ws.Range("D3").Value = "2.231.68"
$ws.Range("E3").Value = "  -3.32%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.621"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.19%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "75.36"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.626"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.33%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0950"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.15"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.103"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.86%  "
$ws.Range("D14").Value = "2.565.84"
$ws.Range("E14").Value = "  -3.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.856"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.12%  "
$ws.Range("D17").Value = "2.235.48"
$ws.Range("E17").Value = "  -2.93%  "
$ws.Range("D18").Value = "41.762.23"
$ws.Range("E18").Value = "  -3.44%  "
$ws.Range("D19").Value = "0.0₃0977"
$ws.Range("E19").Value = "  -2.82%  "
$ws.Range("E20").Value = "  -3.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.80%  "
$ws.Range("E22").Value = "  +2.17%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.46"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -5.97%  "
$ws.Range("E27").Value = "  -5.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.17"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.93%  "
$ws.Range("E29").Value = "  -1.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0833"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.120"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -6.03%  "
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.83"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.89%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0299"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.46%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.91"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.88%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.17"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -8.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "111.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +12.67%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.202"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "60.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.66"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.55%  "
$ws.Range("E46").Value = "  -3.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.64%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -11.37%  "
$ws.Range("E50").Value = "  -2.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.07%  "
